$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the "Overall" (column C) missing-values counts for every
#    "(missing values)" row. The value is simply the sum of the
#    "No IQ" (D) and "Si IQ" (E) missing counts for that row, and must be
#    stored as text (like the surrounding D/E cells) rather than a number.
# ---------------------------------------------------------------------------
$missingRows = @{
    4  = "3"
    6  = "2"
    8  = "3"
    10 = "0"
    14 = "3"
    16 = "5"
    18 = "4"
    20 = "1"
    22 = "1"
    24 = "204"
    26 = "189"
    28 = "204"
    30 = "121"
    32 = "1"
    34 = "195"
    36 = "194"
    38 = "221"
    40 = "3"
    42 = "194"
    44 = "304"
    46 = "1"
    48 = "1"
    50 = "2"
}

foreach ($row in $missingRows.Keys) {
    $cell = $ws.Range("C$row")
    $cell.NumberFormat = "@"
    $cell.Value = $missingRows[$row]
}

# ---------------------------------------------------------------------------
# 2) Several categorical variable labels in column A were shifted up by one
#    row (they now label the first level row of the category instead of the
#    second one).
# ---------------------------------------------------------------------------
$labelMoves = @(
    @{ From = "A52"; To = "A51" },
    @{ From = "A68"; To = "A67" },
    @{ From = "A90"; To = "A89" },
    @{ From = "A94"; To = "A91" }
)

foreach ($move in $labelMoves) {
    $fromRange = $ws.Range($move.From)
    $toRange = $ws.Range($move.To)
    $toRange.Value = $fromRange.Value()
    $fromRange.ClearContents()
}
